$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the header label in B2: was "unnamed: 1_level_1" -> should be "total"
$ws.Range("B2").Value = "total"

# The original data had two stray "label-only" rows that don't belong there
# (row 5 "situação do domicílio" and row 8 "grandes regiões e unidades da
# federação") which pushed every subsequent region's data down by one row
# relative to its label. Removing these two rows corrects the
# label/data alignment for the rest of the table (delete bottom-most first
# so the remaining row index doesn't shift before the second delete).
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
